$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 3.099562666666667
$ws.Cells.Item(2, 8).Value = 9.298688
$ws.Cells.Item(2, 9).Value = 0.2686390288432488
$ws.Cells.Item(2, 10).Value = 0.2686390288432488
$ws.Cells.Item(2, 13).Value = 0.1352566666666667
$ws.Cells.Item(2, 14).Value = 0.40577
$ws.Cells.Item(2, 15).Value = 0.1173241749329269
$ws.Cells.Item(2, 16).Value = 0.1173241749329268
$ws.Cells.Item(2, 17).Value = 0.4192365144177778
$ws.Cells.Item(2, 18).Value = 3.77312862976
$ws.Cells.Item(2, 19).Value = 0.0315178524138169
$ws.Cells.Item(2, 20).Value = 0.0315178524138169

$ws.Cells.Item(3, 7).Value = 3.099562666666667
$ws.Cells.Item(3, 8).Value = 9.298688
$ws.Cells.Item(3, 9).Value = 0.2686390288432488
$ws.Cells.Item(3, 10).Value = 0.2686390288432488
$ws.Cells.Item(3, 15).Value = 0.03951584152489912
$ws.Cells.Item(3, 16).Value = 0.03951584152489912
$ws.Cells.Item(3, 17).Value = 0.1412026436551111
$ws.Cells.Item(3, 18).Value = 1.270823792896
$ws.Cells.Item(3, 19).Value = 0.01061549729117262
$ws.Cells.Item(3, 20).Value = 0.01061549729117262

$ws.Cells.Item(4, 7).Value = 3.099562666666667
$ws.Cells.Item(4, 8).Value = 9.298688
$ws.Cells.Item(4, 9).Value = 0.2686390288432488
$ws.Cells.Item(4, 10).Value = 0.2686390288432488
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.9440163333333332
$ws.Cells.Item(4, 14).Value = 2.832049
$ws.Cells.Item(4, 15).Value = 0.8188575111383802
$ws.Cells.Item(4, 16).Value = 0.8188575111383801
$ws.Cells.Item(4, 17).Value = 2.926037783523555
$ws.Cells.Item(4, 18).Value = 26.334340051712
$ws.Cells.Item(4, 19).Value = 0.2199770865532142
$ws.Cells.Item(4, 20).Value = 0.2199770865532142

$ws.Cells.Item(5, 7).Value = 3.099562666666667
$ws.Cells.Item(5, 8).Value = 9.298688
$ws.Cells.Item(5, 9).Value = 0.2686390288432488
$ws.Cells.Item(5, 10).Value = 0.2686390288432488
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.028017
$ws.Cells.Item(5, 14).Value = 0.084051
$ws.Cells.Item(5, 15).Value = 0.02430247240379386
$ws.Cells.Item(5, 16).Value = 0.02430247240379386
$ws.Cells.Item(5, 17).Value = 0.086840447232
$ws.Cells.Item(5, 18).Value = 0.781564025088
$ws.Cells.Item(5, 19).Value = 0.006528592585045037
$ws.Cells.Item(5, 20).Value = 0.006528592585045036

$ws.Cells.Item(6, 7).Value = 6.189892666666666
$ws.Cells.Item(6, 9).Value = 0.5364778626674904
$ws.Cells.Item(6, 10).Value = 0.5364778626674905
$ws.Cells.Item(6, 13).Value = 0.1352566666666667
$ws.Cells.Item(6, 14).Value = 0.40577
$ws.Cells.Item(6, 15).Value = 0.1173241749329269
$ws.Cells.Item(6, 16).Value = 0.1173241749329268
$ws.Cells.Item(6, 17).Value = 0.8372242491177777
$ws.Cells.Item(6, 18).Value = 7.53501824206
$ws.Cells.Item(6, 19).Value = 0.06294182260724336
$ws.Cells.Item(6, 20).Value = 0.06294182260724336

$ws.Cells.Item(7, 7).Value = 6.189892666666666
$ws.Cells.Item(7, 9).Value = 0.5364778626674904
$ws.Cells.Item(7, 10).Value = 0.5364778626674905
$ws.Cells.Item(7, 15).Value = 0.03951584152489912
$ws.Cells.Item(7, 16).Value = 0.03951584152489912
$ws.Cells.Item(7, 17).Value = 0.2819846870251111
$ws.Cells.Item(7, 18).Value = 2.537862183226
$ws.Cells.Item(7, 19).Value = 0.02119937420278515
$ws.Cells.Item(7, 20).Value = 0.02119937420278515

$ws.Cells.Item(8, 7).Value = 6.189892666666666
$ws.Cells.Item(8, 9).Value = 0.5364778626674904
$ws.Cells.Item(8, 10).Value = 0.5364778626674905
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.9440163333333332
$ws.Cells.Item(8, 14).Value = 2.832049
$ws.Cells.Item(8, 15).Value = 0.8188575111383802
$ws.Cells.Item(8, 16).Value = 0.8188575111383801
$ws.Cells.Item(8, 17).Value = 5.843359778913555
$ws.Cells.Item(8, 18).Value = 52.590238010222
$ws.Cells.Item(8, 19).Value = 0.439298927404739
$ws.Cells.Item(8, 20).Value = 0.439298927404739

$ws.Cells.Item(9, 7).Value = 6.189892666666666
$ws.Cells.Item(9, 9).Value = 0.5364778626674904
$ws.Cells.Item(9, 10).Value = 0.5364778626674905
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.028017
$ws.Cells.Item(9, 14).Value = 0.084051
$ws.Cells.Item(9, 15).Value = 0.02430247240379386
$ws.Cells.Item(9, 16).Value = 0.02430247240379386
$ws.Cells.Item(9, 17).Value = 0.173422222842
$ws.Cells.Item(9, 18).Value = 1.560800005578
$ws.Cells.Item(9, 19).Value = 0.013037738452723
$ws.Cells.Item(9, 20).Value = 0.013037738452723

$ws.Cells.Item(10, 7).Value = 1.888584
$ws.Cells.Item(10, 8).Value = 5.665752
$ws.Cells.Item(10, 9).Value = 0.1636835341659699
$ws.Cells.Item(10, 10).Value = 0.1636835341659699
$ws.Cells.Item(10, 13).Value = 0.1352566666666667
$ws.Cells.Item(10, 14).Value = 0.40577
$ws.Cells.Item(10, 15).Value = 0.1173241749329269
$ws.Cells.Item(10, 16).Value = 0.1173241749329268
$ws.Cells.Item(10, 17).Value = 0.25544357656
$ws.Cells.Item(10, 18).Value = 2.29899218904
$ws.Cells.Item(10, 19).Value = 0.01920403559612796
$ws.Cells.Item(10, 20).Value = 0.01920403559612796

$ws.Cells.Item(11, 7).Value = 1.888584
$ws.Cells.Item(11, 8).Value = 5.665752
$ws.Cells.Item(11, 9).Value = 0.1636835341659699
$ws.Cells.Item(11, 10).Value = 0.1636835341659699
$ws.Cells.Item(11, 15).Value = 0.03951584152489912
$ws.Cells.Item(11, 16).Value = 0.03951584152489912
$ws.Cells.Item(11, 17).Value = 0.086035703176
$ws.Cells.Item(11, 18).Value = 0.7743213285840002
$ws.Cells.Item(11, 19).Value = 0.006468092596337877
$ws.Cells.Item(11, 20).Value = 0.006468092596337878

$ws.Cells.Item(12, 7).Value = 1.888584
$ws.Cells.Item(12, 8).Value = 5.665752
$ws.Cells.Item(12, 9).Value = 0.1636835341659699
$ws.Cells.Item(12, 10).Value = 0.1636835341659699
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 0.9440163333333332
$ws.Cells.Item(12, 14).Value = 2.832049
$ws.Cells.Item(12, 15).Value = 0.8188575111383802
$ws.Cells.Item(12, 16).Value = 0.8188575111383801
$ws.Cells.Item(12, 17).Value = 1.782854142872
$ws.Cells.Item(12, 18).Value = 16.045687285848
$ws.Cells.Item(12, 19).Value = 0.1340334914014801
$ws.Cells.Item(12, 20).Value = 0.1340334914014801

$ws.Cells.Item(13, 7).Value = 1.888584
$ws.Cells.Item(13, 8).Value = 5.665752
$ws.Cells.Item(13, 9).Value = 0.1636835341659699
$ws.Cells.Item(13, 10).Value = 0.1636835341659699
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.028017
$ws.Cells.Item(13, 14).Value = 0.084051
$ws.Cells.Item(13, 15).Value = 0.02430247240379386
$ws.Cells.Item(13, 16).Value = 0.02430247240379386
$ws.Cells.Item(13, 17).Value = 0.052912457928
$ws.Cells.Item(13, 18).Value = 0.476212121352
$ws.Cells.Item(13, 19).Value = 0.003977914572023933
$ws.Cells.Item(13, 20).Value = 0.003977914572023933

$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.3599813333333333
$ws.Cells.Item(14, 8).Value = 1.079944
$ws.Cells.Item(14, 9).Value = 0.03119957432329092
$ws.Cells.Item(14, 10).Value = 0.03119957432329093
$ws.Cells.Item(14, 13).Value = 0.1352566666666667
$ws.Cells.Item(14, 14).Value = 0.40577
$ws.Cells.Item(14, 15).Value = 0.1173241749329269
$ws.Cells.Item(14, 16).Value = 0.1173241749329268
$ws.Cells.Item(14, 17).Value = 0.04868987520888889
$ws.Cells.Item(14, 18).Value = 0.43820887688
$ws.Cells.Item(14, 19).Value = 0.003660464315738637
$ws.Cells.Item(14, 20).Value = 0.003660464315738637

$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.3599813333333333
$ws.Cells.Item(15, 8).Value = 1.079944
$ws.Cells.Item(15, 9).Value = 0.03119957432329092
$ws.Cells.Item(15, 10).Value = 0.03119957432329093
$ws.Cells.Item(15, 15).Value = 0.03951584152489912
$ws.Cells.Item(15, 16).Value = 0.03951584152489912
$ws.Cells.Item(15, 17).Value = 0.01639918962755555
$ws.Cells.Item(15, 18).Value = 0.147592706648
$ws.Cells.Item(15, 19).Value = 0.001232877434603476
$ws.Cells.Item(15, 20).Value = 0.001232877434603476

$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.3599813333333333
$ws.Cells.Item(16, 8).Value = 1.079944
$ws.Cells.Item(16, 9).Value = 0.03119957432329092
$ws.Cells.Item(16, 10).Value = 0.03119957432329093
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.9440163333333332
$ws.Cells.Item(16, 14).Value = 2.832049
$ws.Cells.Item(16, 15).Value = 0.8188575111383802
$ws.Cells.Item(16, 16).Value = 0.8188575111383801
$ws.Cells.Item(16, 17).Value = 0.3398282583617777
$ws.Cells.Item(16, 18).Value = 3.058454325256
$ws.Cells.Item(16, 19).Value = 0.02554800577894692
$ws.Cells.Item(16, 20).Value = 0.02554800577894692

$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 0.6666666666666666
$ws.Cells.Item(17, 7).Value = 0.3599813333333333
$ws.Cells.Item(17, 8).Value = 1.079944
$ws.Cells.Item(17, 9).Value = 0.03119957432329092
$ws.Cells.Item(17, 10).Value = 0.03119957432329093
$ws.Cells.Item(17, 11).Value = 2
$ws.Cells.Item(17, 12).Value = 0.6666666666666666
$ws.Cells.Item(17, 13).Value = 0.028017
$ws.Cells.Item(17, 14).Value = 0.084051
$ws.Cells.Item(17, 15).Value = 0.02430247240379386
$ws.Cells.Item(17, 16).Value = 0.02430247240379386
$ws.Cells.Item(17, 17).Value = 0.010085597016
$ws.Cells.Item(17, 18).Value = 0.09077037314400001
$ws.Cells.Item(17, 19).Value = 0.0007582267940018933
$ws.Cells.Item(17, 20).Value = 0.0007582267940018933

